# The paragraph containing the "<id>" tag for p094r_1 currently holds the
# text "<id>p094r_1</id>" split across three runs (a Courier-New "<id>" run,
# a plain "p094r_1" run, and a Courier-New "</id>" run). The edit merges
# them into a single run carrying the first run's (Courier New / 7f6000 /
# 18pt) formatting, with the full "<id>p094r_1</id>" as its text.
#
# Word's Find/Replace naturally performs this merge: searching across the
# run boundaries for the full logical text and replacing it with identical
# text collapses the matched range into one run that inherits the
# formatting of the first run in the match.

$d = $word.ActiveDocument

$findText = "<id>p094r_1</id>"

$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)

Write-Output "Replaced <id>p094r_1</id>: $found"
